# Apply cryptocurrency price/row updates per the "Updated symbol list" commit.
# All target cells are stored as text (inlineStr) in the workbook, so for the
# "Price" column (D) we force text formatting before assigning numeric-looking
# strings to avoid Excel auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.05"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.00"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.355"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05810"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.377"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.477"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8086"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9212"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1402"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07400"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03192"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03031"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09384"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.853"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001570"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04720"

$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006073"
$ws.Range("E18").Value = "17TigerCashTCH"

$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.001248"
$ws.Range("E19").Value = "18BitKanKAN"

$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004685"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.00008800"
$ws.Range("E21").Value = "20NitroExNTX"

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.593"
$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.150"
$ws.Range("E23").Value = "22BTSETokenBTSE"

$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01079"
$ws.Range("E24").Value = "23OneONEBestin24h"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3180"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1318"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002350"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03849"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006440"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003500"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1065"
$ws.Range("E43").Value = "42BKEXTokenBKK"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009058"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005286"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6855"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001849"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
